# renamed repo, fixed output folder path
# Row 2 (even_MAG-GUT1861.fa) is dropped; every later MAG shifts up by one row.
# New row 2 = old row 3's data (even_MAG-GUT43440.fa), same g__Proteus call.
# New row 3 = a brand-new MAG (even_MAG-GUT49294.fa), predicted/rejected g__Enterobacter.
# New row 4 = old row 7's data (even_MAG-GUT91898.fa), now called g__Proteus.
# Old rows 5 and 6 (even_MAG-GUT49487.fa, even_MAG-GUT91702.fa) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 2: even_MAG-GUT43440.fa (was row 3) ---
$ws.Range("A2").Value = "even_MAG-GUT43440.fa"
$ws.Range("B2").Value = [double]"1.071779194839093e-05"
$ws.Range("C2").Value = [double]"0.0008613058097226647"
$ws.Range("D2").Value = [double]"0.005155611423923531"
$ws.Range("E2").Value = [double]"3.966373574970056e-05"
$ws.Range("F2").Value = [double]"1.118621412034359e-13"
$ws.Range("G2").Value = [double]"0.01414035299303882"
$ws.Range("H2").Value = [double]"0.0003329522403804442"
$ws.Range("I2").Value = [double]"0.00040395809673782"
$ws.Range("J2").Value = [double]"1.251895734106805e-05"
$ws.Range("K2").Value = [double]"0.0003605334811098591"
$ws.Range("L2").Value = [double]"0.0007033738669868567"
$ws.Range("M2").Value = [double]"0.0001991634721122435"
$ws.Range("N2").Value = [double]"8.924440566368064e-05"
$ws.Range("O2").Value = [double]"0.0001027863837946229"
$ws.Range("P2").Value = [double]"9.872290619835539e-05"
$ws.Range("Q2").Value = [double]"5.006167003982559e-08"
$ws.Range("R2").Value = [double]"2.431310147449298e-05"
$ws.Range("S2").Value = [double]"0.7912128980808912"
$ws.Range("T2").Value = [double]"0.1842039773454276"
$ws.Range("U2").Value = [double]"0.00204785584571668"
$ws.Range("V2").Value = [double]"0.7912128980808912"
$ws.Range("W2").Value = "g__Proteus"
$ws.Range("X2").Value = "g__Proteus"

# --- New row 3: even_MAG-GUT49294.fa (new MAG) ---
$ws.Range("A3").Value = "even_MAG-GUT49294.fa"
$ws.Range("B3").Value = [double]"0.01722528475071126"
$ws.Range("C3").Value = [double]"0.01504311103160766"
$ws.Range("D3").Value = [double]"0.03085426517597331"
$ws.Range("E3").Value = [double]"0.00301241496887071"
$ws.Range("F3").Value = [double]"2.150493249622432e-12"
$ws.Range("G3").Value = [double]"0.1330599710450127"
$ws.Range("H3").Value = [double]"0.02549628988327647"
$ws.Range("I3").Value = [double]"0.1015508897392604"
$ws.Range("J3").Value = [double]"0.0007958901134620475"
$ws.Range("K3").Value = [double]"0.1187095437529316"
$ws.Range("L3").Value = [double]"0.03177042534821936"
$ws.Range("M3").Value = [double]"0.04168055727207762"
$ws.Range("N3").Value = [double]"0.06198203807505256"
$ws.Range("O3").Value = [double]"0.09361672620375007"
$ws.Range("P3").Value = [double]"0.0002548747425484968"
$ws.Range("Q3").Value = [double]"4.949150210798939e-05"
$ws.Range("R3").Value = [double]"0.0260971164450182"
$ws.Range("S3").Value = [double]"0.1300556384024548"
$ws.Range("T3").Value = [double]"0.07697024569016575"
$ws.Range("U3").Value = [double]"0.09177522585534865"
$ws.Range("V3").Value = [double]"0.1330599710450127"
$ws.Range("W3").Value = "g__Enterobacter"
$ws.Range("X3").Value = "g__Enterobacter(reject)"

# --- New row 4: even_MAG-GUT91898.fa (was row 7) ---
$ws.Range("A4").Value = "even_MAG-GUT91898.fa"
$ws.Range("B4").Value = [double]"5.666349849918371e-07"
$ws.Range("C4").Value = [double]"0.0002247978637258651"
$ws.Range("D4").Value = [double]"0.0005045264022662864"
$ws.Range("E4").Value = [double]"2.62583190713179e-05"
$ws.Range("F4").Value = [double]"3.104889059675275e-14"
$ws.Range("G4").Value = [double]"0.005438508324598216"
$ws.Range("H4").Value = [double]"4.73552352101178e-05"
$ws.Range("I4").Value = [double]"5.316895259986875e-05"
$ws.Range("J4").Value = [double]"3.938877979820032e-06"
$ws.Range("K4").Value = [double]"0.0001348517097244567"
$ws.Range("L4").Value = [double]"0.0002459404965874328"
$ws.Range("M4").Value = [double]"1.444014565468521e-05"
$ws.Range("N4").Value = [double]"8.81922043862984e-06"
$ws.Range("O4").Value = [double]"1.256687994292176e-05"
$ws.Range("P4").Value = [double]"1.797180384018127e-05"
$ws.Range("Q4").Value = [double]"2.98927850741043e-09"
$ws.Range("R4").Value = [double]"1.628461991320539e-05"
$ws.Range("S4").Value = [double]"0.5523542779706702"
$ws.Range("T4").Value = [double]"0.4405483694665622"
$ws.Range("U4").Value = [double]"0.0003473540869202297"
$ws.Range("V4").Value = [double]"0.5523542779706702"
$ws.Range("W4").Value = "g__Proteus"
$ws.Range("X4").Value = "g__Proteus"

# --- Drop old rows 5 and 6 (and the now-empty trailing row 7) ---
$ws.Range("A5:X7").EntireRow.Delete()
